$d = $word.ActiveDocument

$d.Content.Find.Execute("Horário da Montagem: 01:00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Horário da Montagem: 06:00", 2)

$d.Content.Find.Execute("Desmontagem: 02:00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Desmontagem: 07:00", 2)
